$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (date 44250 -> 44253; Volumen 200 -> 160; N/O/P/S unchanged)
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 160

# Row 3 (date 44252 -> 44250; Volumen 120 -> 200; N 13000->14000; O 14000->15000; P 13500->14500; S 750->806)
$ws.Range("D3").Value = 44250
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806

# Row 4 (date 44253 -> 44252; Volumen 160 -> 120; N 14000->13000; O 15000->14000; P 14500->13500; S 806->750)
$ws.Range("D4").Value = 44252
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13500
$ws.Range("S4").Value = 750
